# sửa nội dung trong luận thiên mã
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "Bạn có thể gặp nhiều điều phát sinh nhưng hầu hết bạn đều xử lý tốt những điều này."
$ws.Range("B10").Value = "Bạn khi ra ngoài gặp nhiều điều may mắn bất ngờ."
$ws.Range("B11").Value = "Bạn gặp nhiều vấn đề trục trặc phát sinh khi di chuyển đi lại, hoặc có nhiều thay đổi trong công việc."

$ws.Range("B11").Select()
